$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.716903209686279
$ws.Range("B1").Value = 2.202841758728027
$ws.Range("C1").Value = 2.342436552047729
$ws.Range("D1").Value = 7.341846466064453
$ws.Range("E1").Value = 0.7762166857719421
